$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remember the existing review comment (currently anchored at C24) so we
#    can re-anchor it one row down (C25) after the row insert below.
# ---------------------------------------------------------------------------
$oldComment = $ws.Comments.Item(1)
$oldCommentText = $oldComment.Text()
$oldComment.Delete()

# ---------------------------------------------------------------------------
# 2. Insert a brand-new row 11 ("Player" / "map_inst_id" / "int" /
#    "当前副本id"), pushing the previous rows 11-38 down to 12-39.
# ---------------------------------------------------------------------------
$ws.Rows.Item(11).Insert()

# ---------------------------------------------------------------------------
# 3. Populate the new row 11 cells.
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "Player"
$ws.Range("B11").Value = "map_inst_id"
$ws.Range("C11").Value = "int"
$ws.Range("D11").Value = "当前副本id"
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0
$ws.Range("P11").Value = 0

# Match the formatting used by the other "int"/"id"-like field rows already
# in the sheet (copy their cell formats onto the new cells).
$ws.Range("C29").Copy()
$ws.Range("C11").PasteSpecial(-4122)

$ws.Range("D20").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("D20").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D20").Copy()
$ws.Range("J11").PasteSpecial(-4122)
$ws.Range("D20").Copy()
$ws.Range("K11").PasteSpecial(-4122)
$ws.Range("D20").Copy()
$ws.Range("P11").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Re-create the comment one row down at C25, preserving its text.
# ---------------------------------------------------------------------------
$ws.Range("C25").AddComment($oldCommentText)

# ---------------------------------------------------------------------------
# 5. Restore the selection shown in the saved sheet (I15).
# ---------------------------------------------------------------------------
$ws.Range("I15").Select()
